$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The relevanceScore column holds percentages as literal text (e.g. "98%"),
# not numeric percent values. Force text format on D2:D4 first so Excel
# doesn't auto-convert the new "100%" strings into numeric percentages.
$ws.Range("D2:D4").NumberFormat = "@"

# Row 2 - Henry Huggins
$ws.Range("A2").Value = "Books for children"
$ws.Range("B2").Value = "Henry Huggins"
$ws.Range("C2").Value = "The book is a classic | It features an average boy whose life is turned upside down when he meets a lovable puppy with a nose for mischief | The book falls under the genre of JUVENILE FICTION"
$ws.Range("D2").Value = "100%"
$ws.Range("E2").Value = "No gap mentioned"

# Row 3 - Anne of Green Gables
$ws.Range("A3").Value = "Books for children"
$ws.Range("B3").Value = "Anne of Green Gables"
$ws.Range("C3").Value = "The book is a classic | It features a talkative eleven-year-old orphan with a heart full of dreams and a desperate longing for a home | The book falls under the genre of JUVENILE FICTION"
$ws.Range("D3").Value = "100%"
$ws.Range("E3").Value = "No gap mentioned"

# Row 4 - The Secret Garden
$ws.Range("A4").Value = "Books for children"
$ws.Range("B4").Value = "The Secret Garden"
$ws.Range("C4").Value = "The book is a classic | It features an orphaned girl who discovers a secret garden and brings it back to life | The book falls under the genre of JUVENILE FICTION"
$ws.Range("D4").Value = "100%"
$ws.Range("E4").Value = "No gap mentioned"
